$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'59.157.60"
$ws.Range("E2").Value = "  +0.53%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'2.583.03"
$ws.Range("E3").Value = "  -0.46%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.01%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'570.06"
$ws.Range("E5").Value = "  +3.42%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'143.51"
$ws.Range("E6").Value = "  +0.07%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.04%  "

# Row 8 - XRP
$ws.Range("D8").Value = "'0.601"
$ws.Range("E8").Value = "  -0.69%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "'2.591.63"
$ws.Range("E9").Value = "  -0.50%  "

# Row 10 - Toncoin
$ws.Range("D10").Value = "'6.67"
$ws.Range("E10").Value = "  -1.97%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +2.89%  "

# Row 12 - TRON
$ws.Range("D12").Value = "'0.156"
$ws.Range("E12").Value = "  +9.87%  "

# Row 13 - Cardano
$ws.Range("D13").Value = "'0.346"
$ws.Range("E13").Value = "  +2.59%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "'3.037.65"
$ws.Range("E14").Value = "  -0.51%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "'59.205.93"
$ws.Range("E15").Value = "  +0.71%  "

# Row 16 - Avalanche
$ws.Range("D16").Value = "'22.46"
$ws.Range("E16").Value = "  +7.65%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  +4.03%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "'2.592.46"
$ws.Range("E18").Value = "  -0.14%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  +1.65%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "'336.50"
$ws.Range("E20").Value = "  -0.02%  "

# Row 21 - Chainlink
$ws.Range("D21").Value = "'10.21"

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +0.41%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.08%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "'64.17"
$ws.Range("E24").Value = "  -3.62%  "

# Row 25 - Polygon
$ws.Range("D25").Value = "'0.452"
$ws.Range("E25").Value = "  +5.46%  "

# Row 26 - Binance-PegBSC-USD
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.40%  "

# Row 27 - Kaspa
$ws.Range("E27").Value = "  +1.61%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("E28").Value = "  +1.75%  "

# Row 29 - PEPE
$ws.Range("D29").Value = "0.0₃0780"
$ws.Range("E29").Value = "  +3.79%  "

# Row 30 - USDe
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  +0.02%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +0.61%  "

# Row 32 - Aptos
$ws.Range("D32").Value = "'6.08"
$ws.Range("E32").Value = "  +2.00%  "

# Row 33 - Monero
$ws.Range("D33").Value = "'157.25"
$ws.Range("E33").Value = "  +2.27%  "

# Row 34 - EthereumClassic
$ws.Range("D34").Value = "'19.07"
$ws.Range("E34").Value = "  +0.73%  "

# Row 35 - NEARProtocol
$ws.Range("D35").Value = "'4.03"
$ws.Range("E35").Value = "  +2.82%  "

# Row 36 - Fetch.AI
$ws.Range("D36").Value = "'0.888"
$ws.Range("E36").Value = "  +7.88%  "

# Row 37 - was SuiNetwork, now ImmutableX (swap with row 38)
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'1.14"
$ws.Range("E37").Value = "  +2.40%  "

# Row 38 - was ImmutableX, now SuiNetwork (swap with row 37)
$ws.Range("B38").Value = "SuiNetwork"
$ws.Range("C38").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D38").Value = "'0.879"
$ws.Range("E38").Value = "  -0.52%  "

# Row 39 - Stacks
$ws.Range("E39").Value = "  +2.35%  "

# Row 40 - OKB
$ws.Range("D40").Value = "'36.88"
$ws.Range("E40").Value = "  -0.08%  "

# Row 41 - Bittensor
$ws.Range("D41").Value = "'293.88"
$ws.Range("E41").Value = "  +4.03%  "

# Row 42 - Filecoin
$ws.Range("E42").Value = "  +1.81%  "

# Row 43 - FirstDigitalUSD
$ws.Range("E43").Value = "  +0.16%  "

# Row 44 - Stellar
$ws.Range("D44").Value = "'0.0977"
$ws.Range("E44").Value = "  +2.05%  "

# Row 45 - Mantle
$ws.Range("D45").Value = "'0.597"
$ws.Range("E45").Value = "  -0.08%  "

# Row 46 - Hedera
$ws.Range("E46").Value = "  +0.91%  "

# Row 47 - EnergySwap
$ws.Range("D47").Value = "'19.19"
$ws.Range("E47").Value = "  +2.72%  "

# Row 48 - WhiteBITCoin
$ws.Range("E48").Value = "  +0.00%  "

# Row 49 - was Aave, now VeChain (swap with row 50)
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0232"
$ws.Range("E49").Value = "  +2.47%  "

# Row 50 - was VeChain, now Aave (swap with row 49)
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'124.13"
$ws.Range("E50").Value = "  +3.22%  "

# Row 51 - InjectiveProtocol
$ws.Range("D51").Value = "'18.53"
$ws.Range("E51").Value = "  +3.91%  "
